$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.87733
$ws.Range("H2").Value = 38.63199
$ws.Range("I2").Value = 0.1584922499374361
$ws.Range("J2").Value = 0.1584922499374361
$ws.Range("M2").Value = 2.733663333333333
$ws.Range("N2").Value = 8.20099
$ws.Range("O2").Value = 0.04037266183309663
$ws.Range("P2").Value = 0.04037266183309663
$ws.Range("Q2").Value = 35.20228485223333
$ws.Range("R2").Value = 316.8205636701
$ws.Range("S2").Value = 0.00639875400989074
$ws.Range("T2").Value = 0.00639875400989074

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.87733
$ws.Range("H3").Value = 38.63199
$ws.Range("I3").Value = 0.1584922499374361
$ws.Range("J3").Value = 0.1584922499374361
$ws.Range("O3").Value = 0.6389522306252696
$ws.Range("P3").Value = 0.6389522306252696
$ws.Range("Q3").Value = 557.1239895557667
$ws.Range("R3").Value = 5014.1159060019
$ws.Range("S3").Value = 0.1012689766343426
$ws.Range("T3").Value = 0.1012689766343426

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.87733
$ws.Range("H4").Value = 38.63199
$ws.Range("I4").Value = 0.1584922499374361
$ws.Range("J4").Value = 0.1584922499374361
$ws.Range("M4").Value = 21.46453166666667
$ws.Range("N4").Value = 64.393595
$ws.Range("O4").Value = 0.3170032929137071
$ws.Range("P4").Value = 0.317003292913707
$ws.Range("Q4").Value = 276.4058575671167
$ws.Range("R4").Value = 2487.65271810405
$ws.Range("S4").Value = 0.05024256513146954
$ws.Range("T4").Value = 0.05024256513146953

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.87733
$ws.Range("H5").Value = 38.63199
$ws.Range("I5").Value = 0.1584922499374361
$ws.Range("J5").Value = 0.1584922499374361
$ws.Range("M5").Value = 0.2486213333333333
$ws.Range("N5").Value = 0.745864
$ws.Range("O5").Value = 0.003671814627926724
$ws.Range("P5").Value = 0.003671814627926724
$ws.Range("Q5").Value = 3.201578954373333
$ws.Range("R5").Value = 28.81421058936
$ws.Range("S5").Value = 0.0005819541617332964
$ws.Range("T5").Value = 0.0005819541617332964

$ws.Range("I6").Value = 0.4359831802722915
$ws.Range("J6").Value = 0.4359831802722916
$ws.Range("M6").Value = 2.733663333333333
$ws.Range("N6").Value = 8.20099
$ws.Range("O6").Value = 0.04037266183309663
$ws.Range("P6").Value = 0.04037266183309663
$ws.Range("Q6").Value = 96.8350446711822
$ws.Range("R6").Value = 871.51540204064
$ws.Range("S6").Value = 0.01760180150205123
$ws.Range("T6").Value = 0.01760180150205123

$ws.Range("I7").Value = 0.4359831802722915
$ws.Range("J7").Value = 0.4359831802722916
$ws.Range("O7").Value = 0.6389522306252696
$ws.Range("P7").Value = 0.6389522306252696
$ws.Range("S7").Value = 0.2785724255500797
$ws.Range("T7").Value = 0.2785724255500798

$ws.Range("I8").Value = 0.4359831802722915
$ws.Range("J8").Value = 0.4359831802722916
$ws.Range("M8").Value = 21.46453166666667
$ws.Range("N8").Value = 64.393595
$ws.Range("O8").Value = 0.3170032929137071
$ws.Range("P8").Value = 0.317003292913707
$ws.Range("Q8").Value = 760.3419402246578
$ws.Range("R8").Value = 6843.077462021921
$ws.Range("S8").Value = 0.1382081038013068
$ws.Range("T8").Value = 0.1382081038013068

$ws.Range("I9").Value = 0.4359831802722915
$ws.Range("J9").Value = 0.4359831802722916
$ws.Range("M9").Value = 0.2486213333333333
$ws.Range("N9").Value = 0.745864
$ws.Range("O9").Value = 0.003671814627926724
$ws.Range("P9").Value = 0.003671814627926724
$ws.Range("Q9").Value = 8.806957911011555
$ws.Range("R9").Value = 79.262621199104
$ws.Range("S9").Value = 0.001600849418853814
$ws.Range("T9").Value = 0.001600849418853814

$ws.Range("G10").Value = 30.51453966666667
$ws.Range("H10").Value = 91.543619
$ws.Range("I10").Value = 0.3755683862706898
$ws.Range("J10").Value = 0.3755683862706898
$ws.Range("M10").Value = 2.733663333333333
$ws.Range("N10").Value = 8.20099
$ws.Range("O10").Value = 0.04037266183309663
$ws.Range("P10").Value = 0.04037266183309663
$ws.Range("Q10").Value = 83.41647822031221
$ws.Range("R10").Value = 750.74830398281
$ws.Range("S10").Value = 0.01516269545410837
$ws.Range("T10").Value = 0.01516269545410837

$ws.Range("G11").Value = 30.51453966666667
$ws.Range("H11").Value = 91.543619
$ws.Range("I11").Value = 0.3755683862706898
$ws.Range("J11").Value = 0.3755683862706898
$ws.Range("O11").Value = 0.6389522306252696
$ws.Range("P11").Value = 0.6389522306252696
$ws.Range("Q11").Value = 1320.179111551154
$ws.Range("R11").Value = 11881.61200396039
$ws.Range("S11").Value = 0.2399702581599901
$ws.Range("T11").Value = 0.2399702581599902

$ws.Range("G12").Value = 30.51453966666667
$ws.Range("H12").Value = 91.543619
$ws.Range("I12").Value = 0.3755683862706898
$ws.Range("J12").Value = 0.3755683862706898
$ws.Range("M12").Value = 21.46453166666667
$ws.Range("N12").Value = 64.393595
$ws.Range("O12").Value = 0.3170032929137071
$ws.Range("P12").Value = 0.317003292913707
$ws.Range("Q12").Value = 654.9803029689228
$ws.Range("R12").Value = 5894.822726720306
$ws.Range("S12").Value = 0.1190564151620958
$ws.Range("T12").Value = 0.1190564151620958

$ws.Range("G13").Value = 30.51453966666667
$ws.Range("H13").Value = 91.543619
$ws.Range("I13").Value = 0.3755683862706898
$ws.Range("J13").Value = 0.3755683862706898
$ws.Range("M13").Value = 0.2486213333333333
$ws.Range("N13").Value = 0.745864
$ws.Range("O13").Value = 0.003671814627926724
$ws.Range("P13").Value = 0.003671814627926724
$ws.Range("Q13").Value = 7.586565537979556
$ws.Range("R13").Value = 68.279089841816
$ws.Range("S13").Value = 0.001379017494495553
$ws.Range("T13").Value = 0.001379017494495553

$ws.Range("G14").Value = 2.433908666666667
$ws.Range("H14").Value = 7.301726
$ws.Range("I14").Value = 0.0299561835195825
$ws.Range("J14").Value = 0.0299561835195825
$ws.Range("M14").Value = 2.733663333333333
$ws.Range("N14").Value = 8.20099
$ws.Range("O14").Value = 0.04037266183309663
$ws.Range("P14").Value = 0.04037266183309663
$ws.Range("Q14").Value = 6.653486878748889
$ws.Range("R14").Value = 59.88138190874
$ws.Range("S14").Value = 0.001209410867046287
$ws.Range("T14").Value = 0.001209410867046287

$ws.Range("G15").Value = 2.433908666666667
$ws.Range("H15").Value = 7.301726
$ws.Range("I15").Value = 0.0299561835195825
$ws.Range("J15").Value = 0.0299561835195825
$ws.Range("O15").Value = 0.6389522306252696
$ws.Range("P15").Value = 0.6389522306252696
$ws.Range("Q15").Value = 105.3004704071178
$ws.Range("R15").Value = 947.70423366406
$ws.Range("S15").Value = 0.01914057028085718
$ws.Range("T15").Value = 0.01914057028085718

$ws.Range("G16").Value = 2.433908666666667
$ws.Range("H16").Value = 7.301726
$ws.Range("I16").Value = 0.0299561835195825
$ws.Range("J16").Value = 0.0299561835195825
$ws.Range("M16").Value = 21.46453166666667
$ws.Range("N16").Value = 64.393595
$ws.Range("O16").Value = 0.3170032929137071
$ws.Range("P16").Value = 0.317003292913707
$ws.Range("Q16").Value = 52.24270964944112
$ws.Range("R16").Value = 470.18438684497
$ws.Range("S16").Value = 0.009496208818834975
$ws.Range("T16").Value = 0.009496208818834973

$ws.Range("G17").Value = 2.433908666666667
$ws.Range("H17").Value = 7.301726
$ws.Range("I17").Value = 0.0299561835195825
$ws.Range("J17").Value = 0.0299561835195825
$ws.Range("M17").Value = 0.2486213333333333
$ws.Range("N17").Value = 0.745864
$ws.Range("O17").Value = 0.003671814627926724
$ws.Range("P17").Value = 0.003671814627926724
$ws.Range("Q17").Value = 0.6051216179182223
$ws.Range("R17").Value = 5.446094561264
$ws.Range("S17").Value = 0.0001099935528440605
$ws.Range("T17").Value = 0.0001099935528440605
